# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# Then turn the data range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) -------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Freeze the header row (row 1) --------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn A1:U80 into a proper Excel Table (ListObject) ------------------
$xlSrcRange = 1   # Microsoft.Office.Interop.Excel.XlListObjectSourceType.xlSrcRange
$xlYes      = 1   # Microsoft.Office.Interop.Excel.XlYesNoGuess.xlYes
$lo = $ws.ListObjects.Add($xlSrcRange, $ws.Range("A1:U80"), $null, $xlYes)
$lo.Name = "Table1"
